$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    8   = @("sd", "Statement-non-opinion")
    19  = @("aa", "Agree/Accept")
    23  = @("sv", "Statement-opinion")
    24  = @("sv", "Statement-opinion")
    27  = @("sd", "Statement-non-opinion")
    33  = @("aa", "Agree/Accept")
    47  = @("%", "Uninterpretable")
    50  = @("sd", "Statement-non-opinion")
    67  = @("aa", "Agree/Accept")
    68  = @("sv", "Statement-opinion")
    69  = @("sv", "Statement-opinion")
    99  = @("sd", "Statement-non-opinion")
    110 = @("aa", "Agree/Accept")
    121 = @("sd", "Statement-non-opinion")
    129 = @("sd", "Statement-non-opinion")
    131 = @("aa", "Agree/Accept")
    134 = @("ba", "Appreciation")
    135 = @("sv", "Statement-opinion")
    139 = @("sv", "Statement-opinion")
    142 = @("sv", "Statement-opinion")
    153 = @("aa", "Agree/Accept")
    162 = @("sd", "Statement-non-opinion")
    173 = @("sd", "Statement-non-opinion")
    184 = @("b", "Acknowledge (Backchannel)")
    213 = @("aa", "Agree/Accept")
    216 = @("b", "Acknowledge (Backchannel)")
    217 = @("sd", "Statement-non-opinion")
    220 = @("sv", "Statement-opinion")
}

foreach ($rowNum in $changes.Keys) {
    $vals = $changes[$rowNum]
    $ws.Range("I$rowNum").Value = $vals[0]
    $ws.Range("J$rowNum").Value = $vals[1]
}
